# ---------------------------------------------------------------------------
# Fix script to python3, automatically remove trailing spaces
#
# This workbook edit:
#   1. Adds a "Variables" sheet (ALMemory key/value scratch area) between
#      "Drinks" and "People".
#   2. Reworks the "People" sheet: replaces the placeholder name list with a
#      real first-name list, adds a "Gender" column driven by two new
#      defined names (female/male) that point at two label cells, and adds
#      a small "Gender list" legend.
#   3. Adds workbook-level defined names "female" and "male".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$locations = $wb.Worksheets.Item("Locations")
$drinks    = $wb.Worksheets.Item("Drinks")
$people    = $wb.Worksheets.Item("People")

# ---------------------------------------------------------------------------
# 1. Defined names used by the new Gender column formulas.
# ---------------------------------------------------------------------------
$wb.Names.Add("female", "=People!`$E`$3")
$wb.Names.Add("male", "=People!`$E`$4")

# ---------------------------------------------------------------------------
# 2. New "Variables" worksheet, inserted right after "Drinks" (i.e. right
#    before "People").
# ---------------------------------------------------------------------------
$vars = $wb.Worksheets.Add($null, $drinks)
$vars.Name = "Variables"

# Re-fetch "People" - its position shifted when the sheet above was
# inserted, and the old reference captured before the insert is stale.
$people = $wb.Worksheets.Item("People")

# Row 1 - headers. Reuse existing header styles from the Locations sheet so
# fonts/borders match the rest of the workbook.
$locations.Range("A1").Copy($vars.Range("A1"))
$vars.Range("A1").Value = "Name"

$locations.Range("B1").Copy($vars.Range("B1"))
$vars.Range("B1").Value = "Location (ALMemory)"

$locations.Range("A1").Copy($vars.Range("C1"))
$vars.Range("C1").Value = "Key 1"
$vars.Range("C1").HorizontalAlignment = -4152

$locations.Range("C1").Copy($vars.Range("D1"))
$vars.Range("D1").Value = "Value 1"

$locations.Range("A1").Copy($vars.Range("E1"))
$vars.Range("E1").Value = "Key 2"
$vars.Range("E1").HorizontalAlignment = -4152

$locations.Range("C1").Copy($vars.Range("F1"))
$vars.Range("F1").Value = "Value 2"

$locations.Range("A1").Copy($vars.Range("G1"))
$vars.Range("G1").Value = "Key 3"
$vars.Range("G1").HorizontalAlignment = -4152

$locations.Range("C1").Copy($vars.Range("H1"))
$vars.Range("H1").Value = "Value 3"

# Row 2 - sample data describing the ALMemory keys written for a guest.
$vars.Range("A2").Value = "John"
$vars.Range("B2").Value = "guest"
$vars.Range("C2").Value = "name"
$vars.Range("C2").HorizontalAlignment = -4152
$vars.Range("D2").Value = "John"
$vars.Range("D2").HorizontalAlignment = -4131
$vars.Range("E2").Value = "drinkId"
$vars.Range("E2").HorizontalAlignment = -4152
$vars.Range("F2").Value = 1
$vars.Range("F2").HorizontalAlignment = -4131

$vars.Range("H5").Select()

# ---------------------------------------------------------------------------
# 3. Rework the "People" sheet.
# ---------------------------------------------------------------------------

# Clear the old placeholder name list (Neo/Bill/Sally/...) before writing
# the new one.
$people.Range("A3:A12").ClearContents()

# Header row additions: a "Gender" column next to the existing "People"
# title, and "name"/"gender" sub-headers.
$people.Range("A1").Copy($people.Range("B1"))
$people.Range("B1").Value = "Gender"

$people.Range("C2").Copy($people.Range("B2"))
$people.Range("B2").Value = "gender"
$people.Range("C2").ClearContents()

$people.Range("E2").Value = "Keys in JSON"
$people.Range("D2").Select()
$locations.Range("E2").Copy($people.Range("D2"))
$people.Range("D2").Value = "Keys in JSON"
$locations.Range("E2").Copy($people.Range("E2"))
$people.Range("E2").Value = "Gender list"

# Labels referenced by the female/male defined names.
$locations.Range("E2").Copy($people.Range("E3"))
$people.Range("E3").Value = "female"
$locations.Range("E2").Copy($people.Range("E4"))
$people.Range("E4").Value = "male"

# New name list - 10 female first names followed by 10 male first names -
# each paired with a formula referencing the matching gender label.
$femaleNames = @("Sophia", "Isabella", "Emma", "Olivia", "Ava", "Emily", "Abigail", "Madison", "Mia", "Chloe")
$maleNames   = @("James", "John", "Robert", "Michael", "William", "David", "Richard", "Charles", "Joseph", "Thomas")

$row = 3
foreach ($name in $femaleNames) {
    $people.Range("A$row").Value = $name
    $people.Range("B$row").Formula = "=female"
    $row++
}
foreach ($name in $maleNames) {
    $people.Range("A$row").Value = $name
    $people.Range("B$row").Formula = "=male"
    $row++
}

# ---------------------------------------------------------------------------
# 4. Selections / active sheet bookkeeping, matching the authored workbook.
# ---------------------------------------------------------------------------
$locations.Activate()
$locations.Range("C3").Select()

$drinks.Activate()
$drinks.Range("A2:C11").Select()

$vars.Activate()
$vars.Range("H5").Select()

$people.Activate()
$people.Range("E4").Select()

Write-Output "Edit complete"
